$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 0.477
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("W2").Value = 0.06317880794701987
$ws.Range("X2").Value = 0.09773871403333599
$ws.Range("Y2").Value = -0.03455990608631612
$ws.Range("AA2").Value = -0.02417302798982188
$ws.Range("AB2").Value = 0.05250516805913231
$ws.Range("AC2").Value = -0.07667819604895421
$ws.Range("AD2").Value = 8.449999999999999
$ws.Range("AF2").Value = 8.449999999999999
$ws.Range("AG2").Value = 8.449999999999999
$ws.Range("AH2").Value = 0.6083513318934485
$ws.Range("AI2").Value = 0.3957845433255269
$ws.Range("AJ2").Value = 0.6083513318934485
$ws.Range("AK2").Value = 0.3957845433255269
$ws.Range("AM2").Value = 0.519
$ws.Range("AO2").Value = -0.7307692307692307
$ws.Range("AQ2").Value = -0.7321772639691715

$ws.Range("B3").Value = "Black Sea Property AS (OB:BSP)"
$ws.Range("K3").Value = 0.477
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("W3").Value = 0.06317880794701987
$ws.Range("X3").Value = 0.09773871403333599
$ws.Range("Y3").Value = -0.03455990608631612
$ws.Range("AA3").Value = -0.02417302798982188
$ws.Range("AB3").Value = 0.05250516805913231
$ws.Range("AC3").Value = -0.07667819604895421
$ws.Range("AD3").Value = 8.449999999999999
$ws.Range("AF3").Value = 8.449999999999999
$ws.Range("AG3").Value = 8.449999999999999
$ws.Range("AH3").Value = 0.6083513318934485
$ws.Range("AI3").Value = 0.3957845433255269
$ws.Range("AJ3").Value = 0.6083513318934485
$ws.Range("AK3").Value = 0.3957845433255269
$ws.Range("AM3").Value = 0.519
$ws.Range("AO3").Value = -0.7307692307692307
$ws.Range("AQ3").Value = -0.7321772639691715
